# Refresh the coin Price / Volume(1h) snapshot in the cryptos worksheet, as
# produced by the scheduled GitHub Actions scrape job. Most rows only get an
# updated Price/Volume; a couple of adjacent coins swap rank (and so swap rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's COM layer auto-coerces a cell's .Value to a Double whenever the
# string parses as a number (e.g. '7.00' -> 7, '612.79' -> 612.78999999...).
# The Price column must stay literal text (it mixes plain decimals with
# thousands-dotted values like '70.569.29'), so numeric-looking prices are
# written with a leading apostrophe - the classic Excel 'force text' quote
# prefix, exactly as if the user typed '7.00 into the cell.
function Set-CellText($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# Row 2 (Bitcoin)
Set-CellText $ws.Range("D2") '70.569.29'
$ws.Range("E2").Value = '  +2.13%  '
# Row 3 (Ethereum)
Set-CellText $ws.Range("D3") '3.561.48'
$ws.Range("E3").Value = '  +1.48%  '
# Row 4 (TetherUSD)
Set-CellText $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.12%  '
# Row 5 (BNB)
Set-CellText $ws.Range("D5") '612.79'
$ws.Range("E5").Value = '  +6.10%  '
# Row 6 (Solana)
Set-CellText $ws.Range("D6") '173.66'
$ws.Range("E6").Value = '  +1.55%  '
# Row 7 (XRP)
Set-CellText $ws.Range("D7") '0.618'
$ws.Range("E7").Value = '  +1.79%  '
# Row 8 (LidoStakedEther)
Set-CellText $ws.Range("D8") '3.557.58'
$ws.Range("E8").Value = '  +1.58%  '
# Row 9 (USDC)
Set-CellText $ws.Range("D9") '0.999'
$ws.Range("E9").Value = '  -0.09%  '
# Row 10 (Dogecoin)
Set-CellText $ws.Range("D10") '0.197'
$ws.Range("E10").Value = '  +4.60%  '
# Row 11 (Toncoin)
Set-CellText $ws.Range("D11") '7.29'
$ws.Range("E11").Value = '  +7.70%  '
# Row 12 (Cardano)
Set-CellText $ws.Range("D12") '0.586'
$ws.Range("E12").Value = '  +0.72%  '
# Row 13 (Avalanche)
Set-CellText $ws.Range("D13") '46.72'
$ws.Range("E13").Value = '  -0.83%  '
# Row 14 (ShibaInu)
$ws.Range("E14").Value = '  +1.77%  '
# Row 15 (WrappedliquidstakedEther2.0)
Set-CellText $ws.Range("D15") '4.134.93'
$ws.Range("E15").Value = '  +1.34%  '
# Row 16 (Polkadot)
Set-CellText $ws.Range("D16") '8.38'
$ws.Range("E16").Value = '  -1.46%  '
# Row 17 (BitcoinCash)
Set-CellText $ws.Range("D17") '616.53'
$ws.Range("E17").Value = '  -1.07%  '
# Row 18 (WrappedEther)
Set-CellText $ws.Range("D18") '3.560.76'
$ws.Range("E18").Value = '  +1.26%  '
# Row 19 (WrappedBTC)
Set-CellText $ws.Range("D19") '70.607.01'
$ws.Range("E19").Value = '  +2.20%  '
# Row 20 (TRON)
$ws.Range("E20").Value = '  -2.04%  '
# Row 21 (Chainlink)
Set-CellText $ws.Range("D21") '17.37'
$ws.Range("E21").Value = '  -0.12%  '
# Row 22 (Polygon)
Set-CellText $ws.Range("D22") '0.882'
$ws.Range("E22").Value = '  -0.14%  '
# Row 23 (Uniswap)
Set-CellText $ws.Range("D23") '9.43'
$ws.Range("E23").Value = '  -15.35%  '
# Row 24 (InternetComputer(DFINITY))
Set-CellText $ws.Range("D24") '15.76'
$ws.Range("E24").Value = '  -0.74%  '
# Row 25 (Litecoin)
Set-CellText $ws.Range("D25") '96.93'
$ws.Range("E25").Value = '  -0.47%  '
# Row 26 (PancakeSwap)
$ws.Range("E26").Value = '  +1.20%  '
# Row 27 (Dai)
$ws.Range("E27").Value = '  +0.09%  '
# Row 28 (ImmutableX)
Set-CellText $ws.Range("D28") '2.62'
$ws.Range("E28").Value = '  -0.47%  '
# Row 29 (EthereumClassic)
Set-CellText $ws.Range("D29") '33.51'
$ws.Range("E29").Value = '  +2.89%  '
# Row 30 (RenderToken)
Set-CellText $ws.Range("D30") '9.05'
$ws.Range("E30").Value = '  -2.80%  '
# Row 31 (Filecoin)
Set-CellText $ws.Range("D31") '8.51'
$ws.Range("E31").Value = '  -0.10%  '
# Row 32 (Stacks)
Set-CellText $ws.Range("D32") '3.06'
$ws.Range("E32").Value = '  -2.76%  '
# Row 33 (NEARProtocol)
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText $ws.Range("D33") '7.00'
$ws.Range("E33").Value = '  -0.04%  '
# Row 34 (Mantle)
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText $ws.Range("D34") '1.30'
$ws.Range("E34").Value = '  -1.45%  '
# Row 35 (Bittensor)
Set-CellText $ws.Range("D35") '574.08'
$ws.Range("E35").Value = '  -9.50%  '
# Row 36 (dogwifhat)
$ws.Range("E36").Value = '  +6.40%  '
# Row 37 (Hedera)
$ws.Range("E37").Value = '  -1.12%  '
# Row 38 (Cosmos)
Set-CellText $ws.Range("D38") '10.83'
$ws.Range("E38").Value = '  +1.01%  '
# Row 39 (OKB)
Set-CellText $ws.Range("D39") '57.28'
$ws.Range("E39").Value = '  +1.04%  '
# Row 40 (VeChain)
$ws.Range("E40").Value = '  +5.32%  '
# Row 41 (FirstDigitalUSD)
Set-CellText $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  +0.11%  '
# Row 42 (Kaspa)
$ws.Range("E42").Value = '  +4.13%  '
# Row 43 (Maker)
Set-CellText $ws.Range("D43") '3.388.71'
$ws.Range("E43").Value = '  +0.24%  '
# Row 44 (TheGraph)
Set-CellText $ws.Range("D44") '0.321'
$ws.Range("E44").Value = '  -1.64%  '
# Row 45 (InjectiveProtocol)
Set-CellText $ws.Range("D45") '33.17'
$ws.Range("E45").Value = '  +1.04%  '
# Row 46 (PEPE)
$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText $ws.Range("D46") '0.0₃0705'
$ws.Range("E46").Value = '  +2.31%  '
# Row 47 (ThetaToken)
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-CellText $ws.Range("D47") '2.96'
$ws.Range("E47").Value = '  +7.91%  '
# Row 48 (Fetch.AI)
Set-CellText $ws.Range("D48") '2.63'
$ws.Range("E48").Value = '  +2.47%  '
# Row 49 (Stellar)
$ws.Range("E49").Value = '  +0.39%  '
# Row 50 (Monero)
Set-CellText $ws.Range("D50") '133.90'
$ws.Range("E50").Value = '  +1.48%  '
# Row 51 (USDe)
$ws.Range("E51").Value = '  -0.02%  '
